$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H92").Value = 557.7
$ws.Range("I92").Value = 660.375
$ws.Range("J92").Value = 147
$ws.Range("K92").Value = 660.375
$ws.Range("L92").Value = 147
$ws.Range("M92").Value = 587.625
$ws.Range("N92").Value = -2643
$ws.Range("H105").Value = 77072
$ws.Range("J105").Value = 77072
$ws.Range("L105").Value = 77072
$ws.Range("N105").Value = -84060
$ws.Range("H107").Value = 171.375
$ws.Range("I107").Value = 178
$ws.Range("K107").Value = 178
$ws.Range("M107").Value = 1742
$ws.Range("H113").Value = 2638.111
$ws.Range("I113").Value = 2237
$ws.Range("J113").Value = 2959
$ws.Range("K113").Value = 2237
$ws.Range("L113").Value = 2959
$ws.Range("M113").Value = 1017
$ws.Range("N113").Value = -9467

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1486.8334
$ws.Range("I2").Value = 678.75
$ws.Range("K2").Value = 678.75
$ws.Range("M2").Value = -565.75
$ws.Range("H32").Value = 3840.4443
$ws.Range("I32").Value = 3840.4443
$ws.Range("K32").Value = 3840.4443
$ws.Range("M32").Value = -3553.4443
$ws.Range("H45").Value = 3964.5
$ws.Range("I45").Value = 2486.2
$ws.Range("K45").Value = 2486.2
$ws.Range("M45").Value = -2109.2
$ws.Range("H102").Value = 2866.5
$ws.Range("I102").Value = 1066.6666
$ws.Range("J102").Value = 4666.3335
$ws.Range("K102").Value = 1066.6666
$ws.Range("L102").Value = 4666.3335
$ws.Range("M102").Value = 555.3334
$ws.Range("N102").Value = -7910.3335
$ws.Range("H108").Value = 105995
$ws.Range("J108").Value = 105995
$ws.Range("L108").Value = 105995
$ws.Range("N108").Value = -113675
$ws.Range("H116").Value = 1486.8334
$ws.Range("I116").Value = 678.75
$ws.Range("K116").Value = 678.75
$ws.Range("M116").Value = 1615.25
$ws.Range("H132").Value = 978.8823
$ws.Range("I132").Value = 978.8823
$ws.Range("K132").Value = 2936.6469
$ws.Range("M132").Value = -406.6468999999997

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1486.8334
$ws.Range("I3").Value = 678.75
$ws.Range("K3").Value = 678.75
$ws.Range("M3").Value = -564.75
$ws.Range("H107").Value = 631.5
$ws.Range("I107").Value = 561
$ws.Range("K107").Value = 561
$ws.Range("M107").Value = 1359
$ws.Range("H134").Value = 5692.1875
$ws.Range("I134").Value = 5621.154
$ws.Range("J134").Value = 6000
$ws.Range("K134").Value = 16863.462
$ws.Range("L134").Value = 18000
$ws.Range("M134").Value = -14328.462
$ws.Range("N134").Value = -23070

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 2278.5
$ws.Range("I7").Value = 1470.7333
$ws.Range("J7").Value = 4701.8
$ws.Range("K7").Value = 1470.7333
$ws.Range("L7").Value = 4701.8
$ws.Range("M7").Value = -1357.7333
$ws.Range("N7").Value = -4927.8
$ws.Range("H17").Value = 1122.5
$ws.Range("I17").Value = 500
$ws.Range("J17").Value = 2990
$ws.Range("K17").Value = 500
$ws.Range("L17").Value = 2990
$ws.Range("M17").Value = -326
$ws.Range("N17").Value = -3338
$ws.Range("H25").Value = 998.5
$ws.Range("I25").Value = 999.5
$ws.Range("J25").Value = 997.5
$ws.Range("K25").Value = 999.5
$ws.Range("L25").Value = 997.5
$ws.Range("M25").Value = -825.5
$ws.Range("N25").Value = -1345.5
$ws.Range("H41").Value = 14321
$ws.Range("J41").Value = 34997
$ws.Range("L41").Value = 34997
$ws.Range("N41").Value = -35853
$ws.Range("H50").Value = 0
$ws.Range("J50").Value = 0
$ws.Range("L50").Value = 0
$ws.Range("N50").ClearContents()
$ws.Range("H51").Value = 24000
$ws.Range("I51").Value = 24000
$ws.Range("K51").Value = 24000
$ws.Range("M51").Value = -23264
$ws.Range("H58").Value = 3640.889
$ws.Range("I58").Value = 1666.5
$ws.Range("K58").Value = 1666.5
$ws.Range("M58").Value = -1463.5
$ws.Range("H60").Value = 18299.5
$ws.Range("J60").Value = 27932.666
$ws.Range("L60").Value = 27932.666
$ws.Range("N60").Value = -28954.666
$ws.Range("H61").Value = 24000
$ws.Range("I61").Value = 24000
$ws.Range("K61").Value = 24000
$ws.Range("M61").Value = -23652
$ws.Range("H105").Value = 1784.7142
$ws.Range("I105").Value = 898
$ws.Range("K105").Value = 898
$ws.Range("M105").Value = 849
$ws.Range("H107").Value = 999
$ws.Range("I107").Value = 999
$ws.Range("K107").Value = 999
$ws.Range("M107").Value = 921
$ws.Range("H136").Value = 3640.889
$ws.Range("I136").Value = 1666.5
$ws.Range("K136").Value = 4999.5
$ws.Range("M136").Value = -2449.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 341.7143
$ws.Range("I2").Value = 499
$ws.Range("J2").Value = 278.8
$ws.Range("K2").Value = 2994
$ws.Range("L2").Value = 1672.8
$ws.Range("M2").Value = -2881
$ws.Range("N2").Value = -1898.8
$ws.Range("H13").Value = 74.888885
$ws.Range("I13").Value = 75.5
$ws.Range("J13").Value = 73.666664
$ws.Range("K13").Value = 226.5
$ws.Range("L13").Value = 220.999992
$ws.Range("M13").Value = -58.5
$ws.Range("N13").Value = -556.999992
$ws.Range("H23").Value = 659.9
$ws.Range("J23").Value = 719.375
$ws.Range("L23").Value = 2158.125
$ws.Range("N23").Value = -2628.125
$ws.Range("H117").Value = 550
$ws.Range("J117").Value = 600
$ws.Range("L117").Value = 1800
$ws.Range("N117").Value = -8684
$ws.Range("H132").Value = 1494.8
$ws.Range("I132").Value = 1699.6666
$ws.Range("J132").Value = 1187.5
$ws.Range("K132").Value = 15296.9994
$ws.Range("L132").Value = 10687.5
$ws.Range("M132").Value = -12766.9994
$ws.Range("N132").Value = -15747.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 8972.4
$ws.Range("I70").Value = 9968
$ws.Range("K70").Value = 9968
$ws.Range("M70").Value = -9698
$ws.Range("H73").Value = 8972.4
$ws.Range("I73").Value = 9968
$ws.Range("K73").Value = 9968
$ws.Range("M73").Value = -9032
$ws.Range("H94").Value = 28332.666
$ws.Range("J94").Value = 28332.666
$ws.Range("L94").Value = 28332.666
$ws.Range("N94").Value = -29684.666
$ws.Range("H113").Value = 799.8
$ws.Range("I113").Value = 733.3333
$ws.Range("J113").Value = 899.5
$ws.Range("K113").Value = 733.3333
$ws.Range("L113").Value = 899.5
$ws.Range("M113").Value = 1436.6667
$ws.Range("N113").Value = -5239.5
$ws.Range("H141").Value = 69999
$ws.Range("I141").Value = 69998
$ws.Range("J141").Value = 70000
$ws.Range("K141").Value = 69998
$ws.Range("L141").Value = 70000
$ws.Range("M141").Value = -64818
$ws.Range("N141").Value = -80360

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1983.75
$ws.Range("I22").Value = 1257.8572
$ws.Range("K22").Value = 1257.8572
$ws.Range("M22").Value = -962.8571999999999
$ws.Range("H27").Value = 1983.75
$ws.Range("I27").Value = 1257.8572
$ws.Range("K27").Value = 1257.8572
$ws.Range("M27").Value = -1150.8572
$ws.Range("H61").Value = 2525.3
$ws.Range("I61").Value = 2303.7144
$ws.Range("J61").Value = 3042.3333
$ws.Range("K61").Value = 2303.7144
$ws.Range("L61").Value = 3042.3333
$ws.Range("M61").Value = -2101.7144
$ws.Range("N61").Value = -3446.3333
$ws.Range("H94").Value = 0
$ws.Range("J94").Value = 0
$ws.Range("L94").Value = 0
$ws.Range("N94").ClearContents()
$ws.Range("H96").Value = 42000
$ws.Range("J96").Value = 42000
$ws.Range("L96").Value = 42000
$ws.Range("N96").Value = -47492
$ws.Range("H113").Value = 2525.3
$ws.Range("I113").Value = 2303.7144
$ws.Range("J113").Value = 3042.3333
$ws.Range("K113").Value = 2303.7144
$ws.Range("L113").Value = 3042.3333
$ws.Range("M113").Value = -133.7143999999998
$ws.Range("N113").Value = -7382.3333

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()
$ws.Range("H94").Value = 49332.332
$ws.Range("J94").Value = 49332.332
$ws.Range("L94").Value = 49332.332
$ws.Range("N94").Value = -51134.332
$ws.Range("H136").Value = 1991.5883
$ws.Range("I136").Value = 1803.625
$ws.Range("J136").Value = 4999
$ws.Range("K136").Value = 5410.875
$ws.Range("L136").Value = 14997
$ws.Range("M136").Value = -2860.875
$ws.Range("N136").Value = -20097
$ws.Range("H140").Value = 79998.5
$ws.Range("I140").Value = 79997
$ws.Range("J140").Value = 80000
$ws.Range("K140").Value = 79997
$ws.Range("L140").Value = 80000
$ws.Range("M140").Value = -74817
$ws.Range("N140").Value = -90360
